$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77: copy style of date cell A76 into A77, then set values/formula
$ws.Range("A76").Copy()
$ws.Range("A77").PasteSpecial(-4122)
$ws.Range("A77").Value = 45645
$ws.Range("B77").Value = 1
$ws.Range("C77").Formula = "=C76+B77"

# Row 78: copy style of date cell A77 into A78, then set values/formula
$ws.Range("A77").Copy()
$ws.Range("A78").PasteSpecial(-4122)
$ws.Range("A78").Value = 45646
$ws.Range("B78").Value = 3
$ws.Range("C78").Formula = "=C77+B78"

$excel.CutCopyMode = 0

$ws.Range("C77").Select()
